# Update cryptos list figures (prices / 1h volume %) as scraped on
# Sat Sep 30 21:46:02 UTC 2023 with GitHub Actions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.064.91"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.678.93"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "215.69"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -3.33%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +1.58%  "

# Row 9 & 10 - Dogecoin and Solana swapped ranking positions
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "21.29"
$ws.Range("E9").Value = "  +5.02%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.0623"
$ws.Range("E10").Value = "  +0.33%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  -0.68%  "

# Row 12 - Wrapped liquid staked Ether 2.0
$ws.Range("D12").Value = "1.918.21"
$ws.Range("E12").Value = "  +0.77%  "

# Row 13 - Wrapped Ether
$ws.Range("D13").Value = "1.689.43"
$ws.Range("E13").Value = "  +1.40%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.70%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.534"
$ws.Range("E15").Value = "  +1.25%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "66.35"

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.070.44"
$ws.Range("E17").Value = "  +0.59%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "8.14"
$ws.Range("E18").Value = "  +2.30%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "235.04"
$ws.Range("E19").Value = "  +0.36%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  +0.33%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.02%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.17%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.24"
$ws.Range("E23").Value = "  +1.24%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -3.43%  "

# Row 25 - Monero
$ws.Range("D25").Value = "147.39"
$ws.Range("E25").Value = "  +0.79%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +1.76%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "16.49"
$ws.Range("E27").Value = "  +3.61%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -1.60%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.14%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.12%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.29%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.35%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.544.50"
$ws.Range("E33").Value = "  +6.17%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.82%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "1.72"
$ws.Range("E35").Value = "  +3.94%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.89%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "0.583"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "0.912"
$ws.Range("E38").Value = "  +1.02%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +2.43%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +7.40%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.04%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "5.53"
$ws.Range("E43").Value = "  -3.60%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -1.02%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.823.25"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46 - TrustWalletToken
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  -0.54%  "

# Row 47 - Quant
$ws.Range("D47").Value = "90.38"
$ws.Range("E47").Value = "  -0.34%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +3.15%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  -0.04%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +1.41%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "8.00"
$ws.Range("E51").Value = "  +6.17%  "
